$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Reposition every shape in the sequence diagram ---
# (slide was dragged/nudged as a whole: dx=-463040 EMU, dy=+685800 EMU)
$newPos = @{
  1 = @(20.540159225463867, 72.0)  # id=81 Rectangle 65: (723900,228600) -> (260860,914400)
  2 = @(33.07913589477539, 96.83039855957031)  # id=2 Rectangle 62: (883145,543946) -> (420105,1229746)
  3 = @(90.3873291015625, 125.46591186523438)  # id=5 Straight Connector 4: (1610959,907617) -> (1147919,1593417)
  4 = @(84.7174072265625, 153.07962036132812)  # id=6 Rectangle 5: (1538951,1258311) -> (1075911,1944111)
  5 = @(234.1848907470703, 87.30882263183594)  # id=16 Rectangle 62: (3437188,423022) -> (2974148,1108822)
  6 = @(282.4840393066406, 125.46591186523438)  # id=17 Straight Connector 16: (4050587,907617) -> (3587547,1593417)
  7 = @(276.8141784667969, 161.54409790039062)  # id=18 Rectangle 17: (3978580,1365810) -> (3515540,2051610)
  8 = @(404.64898681640625, 181.05772399902344)  # id=20 Straight Connector 19: (5602082,1613633) -> (5139042,2299433)
  9 = @(398.64898681640625, 181.05772399902344)  # id=21 Rectangle 20: (5525882,1613633) -> (5062842,2299433)
  10 = @(-3.4598426818847656, 153.37001037597656)  # id=23 Straight Arrow Connector 22: (419100,1261999) -> (-43940,1947799)
  11 = @(-33.459842681884766, 132.0)  # id=26 TextBox 25: (38100,990600) -> (-424940,1676400)
  12 = @(289.2072448730469, 173.08189392089844)  # id=28 Straight Arrow Connector 27: (4135972,1512340) -> (3672932,2198140)
  13 = @(218.92764282226562, 249.5961456298828)  # id=29 TextBox 28: (3243421,2484071) -> (2780381,3169871)
  14 = @(287.09197998046875, 201.8922882080078)  # id=34 Straight Arrow Connector 33: (4109108,1878232) -> (3646068,2564032)
  15 = @(96.7174072265625, 222.0)  # id=35 Straight Arrow Connector 34: (1691351,2133600) -> (1228311,2819400)
  16 = @(-6.459921360015869, 384.0)  # id=37 Straight Arrow Connector 36: (380999,4191000) -> (-82041,4876800)
  17 = @(398.69671630859375, 269.06402587890625)  # id=65 Rectangle 64: (5526488,2731313) -> (5063448,3417113)
  18 = @(412.3852844238281, 338.2362365722656)  # id=75 Straight Arrow Connector 74: (5700333,3609800) -> (5237293,4295600)
  19 = @(418.8734130859375, 271.4012756347656)  # id=79 TextBox 78: (5782732,2760996) -> (5319692,3446796)
  20 = @(111.98023986816406, 141.0984344482422)  # id=80 TextBox 79: (1885189,1106150) -> (1422149,1791950)
  21 = @(221.23741149902344, 352.50994873046875)  # id=82 TextBox 81: (3272755,3791076) -> (2809715,4476876)
  22 = @(14.348819732666016, 364.70086669921875)  # id=83 TextBox 82: (645270,3945901) -> (182230,4631701)
  23 = @(579.8384399414062, 273.063720703125)  # id=84 Rectangle 62: (7826988,2782109) -> (7363948,3467909)
  24 = @(633.9501953125, 337.40386962890625)  # id=86 Rectangle 85: (8514207,3599229) -> (8051167,4285029)
  25 = @(411.2374267578125, 361.0697021484375)  # id=90 Straight Arrow Connector 89: (5685755,3899785) -> (5222715,4585585)
  26 = @(178.09071350097656, 204.0624542236328)  # id=93 TextBox 92: (2724792,1905793) -> (2261752,2591593)
  27 = @(440.6708679199219, 229.6177215576172)  # id=40 Rectangle 62: (6059560,2230345) -> (5596520,2916145)
  28 = @(513.9691772460938, 262.9217529296875)  # id=46 Straight Connector 45: (6990448,2653306) -> (6527408,3339106)
  29 = @(505.86505126953125, 286.92181396484375)  # id=49 Rectangle 48: (6887526,2958107) -> (6424486,3643907)
  30 = @(411.2374267578125, 288.2790832519531)  # id=50 Straight Arrow Connector 49: (5685755,2975344) -> (5222715,3661144)
  31 = @(394.4756164550781, 389.2922058105469)  # id=4 TextBox 3: (5472880,4258211) -> (5009840,4944011)
  32 = @(360.07373046875, 153.23370361328125)  # id=19 Rectangle 62: (5035976,1260268) -> (4572936,1946068)
  33 = @(96.7174072265625, 269.0641174316406)  # id=66 Straight Arrow Connector 65: (1691351,2731314) -> (1228311,3417114)
  34 = @(98.04763793945312, 161.39512634277344)  # id=67 Straight Arrow Connector 66: (1708245,1363918) -> (1245205,2049718)
  35 = @(96.76834869384766, 371.8316650390625)  # id=69 Straight Arrow Connector 68: (1691998,4036462) -> (1228958,4722262)
  36 = @(414.8223876953125, 314.0008850097656)  # id=74 Straight Arrow Connector 73: (5731284,3302011) -> (5268244,3987811)
  37 = @(433.46875, 298.1287536621094)  # id=85 TextBox 84: (5968093,3100435) -> (5505053,3786235)
  38 = @(441.6382141113281, 323.69561767578125)  # id=88 TextBox 87: (6071845,3425134) -> (5608805,4110934)
  39 = @(639.9501953125, 305.3990783691406)  # id=89 Straight Connector 88: (8590407,3192768) -> (8127367,3878568)
}
foreach ($idx in $newPos.Keys) {
  $pos = $newPos[$idx]
  $sh = $s.Shapes.Item($idx)
  $sh.Left = $pos[0]
  $sh.Top = $pos[1]
}

# --- Refresh the auto "datetimeFigureOut" footer placeholders ---
$newDate = "10/22/2017"
foreach ($sh in $p.SlideMaster.Shapes) {
  if ($sh.Name -like "Date Placeholder*") {
    $sh.TextFrame.TextRange.Text = $newDate
  }
}
foreach ($cl in $p.SlideMaster.CustomLayouts) {
  foreach ($sh in $cl.Shapes) {
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = $newDate
    }
  }
}
foreach ($sh in $p.NotesMaster.Shapes) {
  if ($sh.Name -like "Date Placeholder*") {
    $sh.TextFrame.TextRange.Text = $newDate
  }
}

